$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows above the old row 9 (the "Tuesday" entry). This
#    shifts old row 9 -> 11, old row 10 -> 12, old rows 11/12 -> 13/14, and
#    so on down to the bottom of the sheet. Row 8 (previously blank) keeps
#    its row number.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# ---------------------------------------------------------------------------
# 2. Copy cell formatting onto the newly created rows so the new data picks
#    up the same look as the rest of the table. Every copy/paste below uses
#    single-cell (1x1) ranges so the pasted format can never spill into
#    neighbouring columns.
#    Row 12 (a normal, unbordered data row) supplies the format for the new
#    rows 8 and 9, as well as for the right-hand (F/G/H/J) cells of row 10/11.
#    Row 13/14 (a bordered, merged two-entry day) supplies the format for
#    the B/C/D columns of the new two-entry day in rows 10/11.
# ---------------------------------------------------------------------------
foreach ($col in @("B", "C", "D", "F", "G", "H", "J")) {
    $ws.Range("$col`12").Copy()
    $ws.Range("$col`8").PasteSpecial(-4122)
    $ws.Range("$col`12").Copy()
    $ws.Range("$col`9").PasteSpecial(-4122)
}

foreach ($col in @("F", "G", "H", "J")) {
    $ws.Range("$col`12").Copy()
    $ws.Range("$col`10").PasteSpecial(-4122)
    $ws.Range("$col`12").Copy()
    $ws.Range("$col`11").PasteSpecial(-4122)
}

foreach ($col in @("B", "C", "D")) {
    $ws.Range("$col`13").Copy()
    $ws.Range("$col`10").PasteSpecial(-4122)
    $ws.Range("$col`14").Copy()
    $ws.Range("$col`11").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the new / changed values.
#    Shared-string cells are written in the same order the author apparently
#    typed them in, so new strings land at the expected shared-string index:
#      1st new string -> J10 (Amberdata logo icons)
#      2nd new string -> D9  (Wednesday)
#      3rd new string -> J9  (social media header graphics)
#      4th new string -> J8  (continued editing social media headers)
# ---------------------------------------------------------------------------

# Row 10: new second task on Tuesday 2018-07-24 (#3)
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 43305
$ws.Range("D10").Value = "Tuesday"
$ws.Range("F10").Value = 0.79166666666666663
$ws.Range("G10").Value = 0.83333333333333337
$ws.Range("H10").Value = 1
$ws.Range("J10").Value = "▫ Created PNG/JPG Amberdata Logo icons in 8 different sizes"

# Row 9: Wednesday 2018-07-25 (#4)
$ws.Range("C9").Value = 43306
$ws.Range("D9").Value = "Wednesday"
$ws.Range("F9").Value = 0.45833333333333331
$ws.Range("G9").Value = 0.54166666666666663
$ws.Range("J9").Value = "▫ Created graphics for social media headers"

# Row 8: Thursday 2018-07-26 (#5)
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 43307
$ws.Range("D8").Value = "Thursday"
$ws.Range("F8").Value = 0.54166666666666663
$ws.Range("G8").Value = 0.66666666666666663
$ws.Range("H8").Value = 3
$ws.Range("J8").Value = "▫ Contiued editing graphics for social media headers"

# Row 9 unchanged values (kept identical to before the edit)
$ws.Range("B9").Value = 4
$ws.Range("H9").Value = 2

# Row 11: continuation of the Tuesday block - keep the original task that
# used to live in row 9, clear the now-redundant date/day columns (they will
# be merged with row 10) and keep the time/hours/task values.
$ws.Range("F11").Value = 0.58333333333333337
$ws.Range("G11").Value = 0.66666666666666663
$ws.Range("H11").Value = 2
$ws.Range("J11").Value = "▫ Created Aion Announcement Header"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""

# ---------------------------------------------------------------------------
# 4. Merge the B/C/D columns of the new two-entry Tuesday block, matching the
#    pattern already used by the other multi-entry day further down.
# ---------------------------------------------------------------------------
$ws.Range("B10:B11").Merge()
$ws.Range("C10:C11").Merge()
$ws.Range("D10:D11").Merge()

# ---------------------------------------------------------------------------
# 5. Two new blank (but date-formatted) rows at the bottom of the log.
# ---------------------------------------------------------------------------
$ws.Range("C51").Copy()
$ws.Range("C52").PasteSpecial(-4122)
$ws.Range("C51").Copy()
$ws.Range("C53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Restore the selected cell shown in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("J22").Select()
